$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = 43647
$ws.Range("M1").NumberFormat = "mm-dd-yy"
$ws.Range("M1").Interior.Color = 49407

Write-Host "done"
